$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1, matching the style used by the other header cells (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Numeric data for I2:J39
$data = @{
    2 = @(6, 7)
    3 = @(7, 7)
    4 = @(8, 8)
    5 = @(9, 9)
    6 = @(7, 7)
    7 = @(8, 8)
    8 = @(9, 9)
    9 = @(6, 7)
    10 = @(10, 11)
    11 = @(6, 6)
    12 = @(2, 3)
    13 = @(9, 9)
    14 = @(7, 7)
    15 = @(9, 9)
    16 = @(7, 7)
    17 = @(9, 9)
    18 = @(6, 7)
    19 = @(9, 9)
    20 = @(9, 9)
    21 = @(8, 8)
    22 = @(5, 5)
    23 = @(6, 6)
    24 = @(5, 5)
    25 = @(7, 7)
    26 = @(8, 8)
    27 = @(8, 8)
    28 = @(10, 10)
    29 = @(7, 7)
    30 = @(4, 4)
    31 = @(7, 7)
    32 = @(5, 5)
    33 = @(6, 6)
    34 = @(8, 8)
    35 = @(7, 7)
    36 = @(7, 7)
    37 = @(4, 4)
    38 = @(5, 5)
    39 = @(4, 4)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}
